# Add new power plant / generation technology rows to the Electricity Source
# subscript on the "PDiCCpDoC" sheet (issues #280 and #99).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("PDiCCpDoC")

# The existing data rows (B2:B18) previously carried an explicit "applied
# number format" style that did nothing visually (numFmtId 0 == General).
# Strip that formatting so those cells fall back to the default style,
# matching how the sheet looks after the new rows were typed in directly.
$ws.Range("B2:B18").ClearFormats()

# New Electricity Source categories being added to the subscript.
$newSources = @(
    "hard coal w CCS",
    "natural gas combined cycle w CCS",
    "biomass w CCS",
    "lignite w CCS",
    "small modular reactor",
    "hydrogen"
)

$row = 19
foreach ($name in $newSources) {
    $ws.Cells.Item($row, 1).Value = $name
    $ws.Cells.Item($row, 2).Value = 0
    $ws.Cells.Item($row, 2).NumberFormat = "0"
    $row = $row + 1
}

# Leave the cursor where the author left it (just past the new data) and
# restore "About" as the active sheet/tab, matching the saved selection
# state in the workbook.
$ws.Range("A25").Select() | Out-Null
$wsAbout = $wb.Worksheets.Item("About")
$wsAbout.Activate() | Out-Null
